# Update cached TPM-derived values on the active (only) worksheet.
# Each row of the LR-pairs table is recomputed from the new TPM input;
# write the refreshed numbers straight into the cells the diff touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 1.10087;                "H2" = 3.30261
    "I2" = 0.1843884439613191;     "J2" = 0.1843884439613191
    "M2" = 5.916202333333334;      "N2" = 17.748607
    "O2" = 0.3515586392055965;     "P2" = 0.3515586392055965
    "Q2" = 6.512969662696667;      "R2" = 58.61672696427
    "S2" = 0.06482335044427873;    "T2" = 0.06482335044427873

    "G3" = 1.10087;                "H3" = 3.30261
    "I3" = 0.1843884439613191;     "J3" = 0.1843884439613191
    "O3" = 0.6159539016771971;     "P3" = 0.6159539016771971
    "Q3" = 11.41115201807667;      "R3" = 102.70036816269
    "S3" = 0.1135747814821617;     "T3" = 0.1135747814821617

    "G4" = 1.10087;                "H4" = 3.30261
    "I4" = 0.1843884439613191;     "J4" = 0.1843884439613191
    "M4" = 0.5467150000000001
    "O4" = 0.03248745911720639;    "P4" = 0.03248745911720639
    "Q4" = 0.6018621420500001;     "R4" = 5.416759278450001
    "S4" = 0.005990312034878655;   "T4" = 0.005990312034878655

    "I5" = 0.4503925067925547;     "J5" = 0.4503925067925547
    "M5" = 5.916202333333334;      "N5" = 17.748607
    "O5" = 0.3515586392055965;     "P5" = 0.3515586392055965
    "Q5" = 15.90876667770556;      "R5" = 143.17890009935
    "S5" = 0.1583393767963879;     "T5" = 0.1583393767963879

    "I6" = 0.4503925067925547;     "J6" = 0.4503925067925547
    "O6" = 0.6159539016771971;     "P6" = 0.6159539016771971
    "S6" = 0.2774210218450476;     "T6" = 0.2774210218450475

    "I7" = 0.4503925067925547;     "J7" = 0.4503925067925547
    "M7" = 0.5467150000000001
    "O7" = 0.03248745911720639;    "P7" = 0.03248745911720639
    "S7" = 0.01463210815111922;    "T7" = 0.01463210815111922

    "H8" = 6.541494999999999
    "I8" = 0.3652190492461261;     "J8" = 0.3652190492461262
    "M8" = 5.916202333333334;      "N8" = 17.748607
    "O8" = 0.3515586392055965;     "P8" = 0.3515586392055965
    "Q8" = 12.90026932749611;      "R8" = 116.102423947465
    "S8" = 0.1283959119649299;     "T8" = 0.1283959119649299

    "H9" = 6.541494999999999
    "I9" = 0.3652190492461261;     "J9" = 0.3652190492461262
    "O9" = 0.6159539016771971;     "P9" = 0.6159539016771971
    "S9" = 0.2249580983499878;     "T9" = 0.2249580983499878

    "H10" = 6.541494999999999
    "I10" = 0.3652190492461261;    "J10" = 0.3652190492461262
    "M10" = 0.5467150000000001
    "O10" = 0.03248745911720639;   "P10" = 0.03248745911720639
    "S10" = 0.01186503893120851;   "T10" = 0.01186503893120851
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
